$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format so numeric-looking strings
# (e.g. "1.010", "0.00000000349") are preserved exactly as text,
# matching the original inlineStr cell contents.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.356.04"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "2.102.76"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "344.28"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.5222"
$ws.Range("E7").Value = "  +1.70%  "
$ws.Range("D8").Value = "0.4422"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "54.77"
$ws.Range("E9").Value = "  +2.95%  "
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").Value = "1.170"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "24.79"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "8.646"
$ws.Range("E13").Value = "  +6.14%  "
$ws.Range("D14").Value = "2.143.35"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").Value = "6.912"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "101.62"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "0.00001159"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "21.14"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = "6.369"
$ws.Range("E21").Value = "  +3.13%  "
$ws.Range("D22").Value = "1.006"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "30.416.74"
$ws.Range("E23").Value = "  +2.29%  "
$ws.Range("D24").Value = "12.54"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "2.300"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "21.88"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "162.45"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "2.512"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").Value = "133.46"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").Value = "1.134"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1051"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").Value = "1.664"
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "6.729"
$ws.Range("E33").Value = "  +10.98%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "6.221"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").Value = "3.921"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").Value = "10.30"
$ws.Range("D37").Value = "0.02626"
$ws.Range("E37").Value = "  +2.49%  "
$ws.Range("D38").Value = "0.06768"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").Value = "0.7013"
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.344"
$ws.Range("E40").Value = "  +3.88%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "12.51"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").Value = "0.2221"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").Value = "0.6822"
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("D44").Value = "14.49"
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("D45").Value = "2.348"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "1.398"
$ws.Range("E47").Value = "  +20.51%  "
$ws.Range("D48").Value = "3.640"
$ws.Range("D49").Value = "0.00000000349"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("D50").Value = "1.212"
$ws.Range("E50").Value = "  +9.08%  "
$ws.Range("D51").Value = "1.217"
$ws.Range("E51").Value = "  -0.06%  "
